$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -191335.61
$ws.Range("D2").Value = -198285.63
$ws.Range("E2").Value = -182703.46
$ws.Range("F2").Value = -206749.08
$ws.Range("G2").Value = -174255.89
$ws.Range("H2").Value = -953329.67

$ws.Range("C3").Value = -54488.72
$ws.Range("D3").Value = -46939.96
$ws.Range("E3").Value = -54448.97
$ws.Range("F3").Value = -58853.76
$ws.Range("G3").Value = -47277.03
$ws.Range("H3").Value = -262008.44

$ws.Range("C4").Value = 591845.34
$ws.Range("D4").Value = 585293.16
$ws.Range("E4").Value = 643798.41
$ws.Range("F4").Value = 537409.59
$ws.Range("G4").Value = 647993.1
$ws.Range("H4").Value = 3006339.6

$ws.Range("C6").Value = 346021.01
$ws.Range("D6").Value = 340067.57
$ws.Range("E6").Value = 406645.98
$ws.Range("F6").Value = 271806.75
$ws.Range("G6").Value = 426460.18
$ws.Range("H6").Value = 1791001.49
